$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cell C3 to a single space, mirroring how the author's fix for
# "empty cells" handling adds a whitespace-only value to the test data.
$ws.Range("C3").Value = " "

# Move the active selection to B6 (reflecting the final selection state
# in the saved workbook).
$ws.Range("B6").Select()
